# Auto-generated: append updated COVID case data rows through 2021-12-08 (8/12)
# as described by commit message 'aggiornamento fino a 8/12'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 386
$lastRow = 464
$numRows = 79

# Propagate the date-style formatting (style index 2, i.e. YYYY-MM-DD HH:MM:SS,
# centered thin-bordered cell) from column A of the last existing row down through
# column A of all the newly appended rows, matching the existing data pattern.
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = New-Object "object[,]" $numRows,4
$data[0,0] = 44460; $data[0,1] = 0; $data[0,2] = 3; $data[0,3] = 36.46086533787069
$data[1,0] = 44461; $data[1,1] = 0; $data[1,2] = 3; $data[1,3] = 36.46086533787069
$data[2,0] = 44462; $data[2,1] = 1; $data[2,2] = 3; $data[2,3] = 36.46086533787069
$data[3,0] = 44463; $data[3,1] = 0; $data[3,2] = 2; $data[3,3] = 24.30724355858046
$data[4,0] = 44464; $data[4,1] = 2; $data[4,2] = 4; $data[4,3] = 48.61448711716091
$data[5,0] = 44465; $data[5,1] = 0; $data[5,2] = 4; $data[5,3] = 48.61448711716091
$data[6,0] = 44466; $data[6,1] = 0; $data[6,2] = 3; $data[6,3] = 36.46086533787069
$data[7,0] = 44467; $data[7,1] = 0; $data[7,2] = 3; $data[7,3] = 36.46086533787069
$data[8,0] = 44468; $data[8,1] = 2; $data[8,2] = 5; $data[8,3] = 60.76810889645115
$data[9,0] = 44469; $data[9,1] = 0; $data[9,2] = 4; $data[9,3] = 48.61448711716091
$data[10,0] = 44470; $data[10,1] = 2; $data[10,2] = 6; $data[10,3] = 72.92173067574137
$data[11,0] = 44471; $data[11,1] = 0; $data[11,2] = 4; $data[11,3] = 48.61448711716091
$data[12,0] = 44472; $data[12,1] = 2; $data[12,2] = 6; $data[12,3] = 72.92173067574137
$data[13,0] = 44473; $data[13,1] = 0; $data[13,2] = 6; $data[13,3] = 72.92173067574137
$data[14,0] = 44474; $data[14,1] = 0; $data[14,2] = 6; $data[14,3] = 72.92173067574137
$data[15,0] = 44475; $data[15,1] = 0; $data[15,2] = 4; $data[15,3] = 48.61448711716091
$data[16,0] = 44476; $data[16,1] = 0; $data[16,2] = 4; $data[16,3] = 48.61448711716091
$data[17,0] = 44477; $data[17,1] = 2; $data[17,2] = 4; $data[17,3] = 48.61448711716091
$data[18,0] = 44478; $data[18,1] = 0; $data[18,2] = 4; $data[18,3] = 48.61448711716091
$data[19,0] = 44479; $data[19,1] = 0; $data[19,2] = 2; $data[19,3] = 24.30724355858046
$data[20,0] = 44480; $data[20,1] = 0; $data[20,2] = 2; $data[20,3] = 24.30724355858046
$data[21,0] = 44481; $data[21,1] = 0; $data[21,2] = 2; $data[21,3] = 24.30724355858046
$data[22,0] = 44482; $data[22,1] = 0; $data[22,2] = 2; $data[22,3] = 24.30724355858046
$data[23,0] = 44483; $data[23,1] = 0; $data[23,2] = 2; $data[23,3] = 24.30724355858046
$data[24,0] = 44484; $data[24,1] = 0; $data[24,2] = 0; $data[24,3] = 0.0
$data[25,0] = 44485; $data[25,1] = 0; $data[25,2] = 0; $data[25,3] = 0.0
$data[26,0] = 44486; $data[26,1] = 0; $data[26,2] = 0; $data[26,3] = 0.0
$data[27,0] = 44487; $data[27,1] = 0; $data[27,2] = 0; $data[27,3] = 0.0
$data[28,0] = 44488; $data[28,1] = 0; $data[28,2] = 0; $data[28,3] = 0.0
$data[29,0] = 44489; $data[29,1] = 0; $data[29,2] = 0; $data[29,3] = 0.0
$data[30,0] = 44490; $data[30,1] = 0; $data[30,2] = 0; $data[30,3] = 0.0
$data[31,0] = 44491; $data[31,1] = 1; $data[31,2] = 1; $data[31,3] = 12.15362177929023
$data[32,0] = 44492; $data[32,1] = 0; $data[32,2] = 1; $data[32,3] = 12.15362177929023
$data[33,0] = 44493; $data[33,1] = 1; $data[33,2] = 2; $data[33,3] = 24.30724355858046
$data[34,0] = 44494; $data[34,1] = 0; $data[34,2] = 2; $data[34,3] = 24.30724355858046
$data[35,0] = 44495; $data[35,1] = 0; $data[35,2] = 2; $data[35,3] = 24.30724355858046
$data[36,0] = 44496; $data[36,1] = 1; $data[36,2] = 3; $data[36,3] = 36.46086533787069
$data[37,0] = 44497; $data[37,1] = 1; $data[37,2] = 4; $data[37,3] = 48.61448711716091
$data[38,0] = 44498; $data[38,1] = 1; $data[38,2] = 4; $data[38,3] = 48.61448711716091
$data[39,0] = 44499; $data[39,1] = 0; $data[39,2] = 4; $data[39,3] = 48.61448711716091
$data[40,0] = 44500; $data[40,1] = 0; $data[40,2] = 3; $data[40,3] = 36.46086533787069
$data[41,0] = 44501; $data[41,1] = 2; $data[41,2] = 5; $data[41,3] = 60.76810889645115
$data[42,0] = 44502; $data[42,1] = 3; $data[42,2] = 8; $data[42,3] = 97.22897423432183
$data[43,0] = 44503; $data[43,1] = 0; $data[43,2] = 7; $data[43,3] = 85.0753524550316
$data[44,0] = 44504; $data[44,1] = 0; $data[44,2] = 6; $data[44,3] = 72.92173067574137
$data[45,0] = 44505; $data[45,1] = 2; $data[45,2] = 7; $data[45,3] = 85.0753524550316
$data[46,0] = 44506; $data[46,1] = 0; $data[46,2] = 7; $data[46,3] = 85.0753524550316
$data[47,0] = 44507; $data[47,1] = 2; $data[47,2] = 9; $data[47,3] = 109.3825960136121
$data[48,0] = 44508; $data[48,1] = 1; $data[48,2] = 8; $data[48,3] = 97.22897423432183
$data[49,0] = 44509; $data[49,1] = 0; $data[49,2] = 5; $data[49,3] = 60.76810889645115
$data[50,0] = 44510; $data[50,1] = 0; $data[50,2] = 5; $data[50,3] = 60.76810889645115
$data[51,0] = 44511; $data[51,1] = 0; $data[51,2] = 5; $data[51,3] = 60.76810889645115
$data[52,0] = 44512; $data[52,1] = 0; $data[52,2] = 3; $data[52,3] = 36.46086533787069
$data[53,0] = 44513; $data[53,1] = 0; $data[53,2] = 3; $data[53,3] = 36.46086533787069
$data[54,0] = 44514; $data[54,1] = 1; $data[54,2] = 2; $data[54,3] = 24.30724355858046
$data[55,0] = 44515; $data[55,1] = 0; $data[55,2] = 1; $data[55,3] = 12.15362177929023
$data[56,0] = 44516; $data[56,1] = 3; $data[56,2] = 4; $data[56,3] = 48.61448711716091
$data[57,0] = 44517; $data[57,1] = 0; $data[57,2] = 4; $data[57,3] = 48.61448711716091
$data[58,0] = 44518; $data[58,1] = 3; $data[58,2] = 7; $data[58,3] = 85.0753524550316
$data[59,0] = 44519; $data[59,1] = 0; $data[59,2] = 7; $data[59,3] = 85.0753524550316
$data[60,0] = 44520; $data[60,1] = 1; $data[60,2] = 8; $data[60,3] = 97.22897423432183
$data[61,0] = 44521; $data[61,1] = 1; $data[61,2] = 8; $data[61,3] = 97.22897423432183
$data[62,0] = 44522; $data[62,1] = 4; $data[62,2] = 12; $data[62,3] = 145.8434613514827
$data[63,0] = 44523; $data[63,1] = 0; $data[63,2] = 9; $data[63,3] = 109.3825960136121
$data[64,0] = 44524; $data[64,1] = 3; $data[64,2] = 12; $data[64,3] = 145.8434613514827
$data[65,0] = 44525; $data[65,1] = 1; $data[65,2] = 10; $data[65,3] = 121.5362177929023
$data[66,0] = 44526; $data[66,1] = 8; $data[66,2] = 18; $data[66,3] = 218.7651920272241
$data[67,0] = 44527; $data[67,1] = 0; $data[67,2] = 17; $data[67,3] = 206.6115702479339
$data[68,0] = 44528; $data[68,1] = 5; $data[68,2] = 21; $data[68,3] = 255.2260573650948
$data[69,0] = 44529; $data[69,1] = 1; $data[69,2] = 18; $data[69,3] = 218.7651920272241
$data[70,0] = 44530; $data[70,1] = 4; $data[70,2] = 22; $data[70,3] = 267.379679144385
$data[71,0] = 44531; $data[71,1] = 1; $data[71,2] = 20; $data[71,3] = 243.0724355858046
$data[72,0] = 44532; $data[72,1] = 1; $data[72,2] = 20; $data[72,3] = 243.0724355858046
$data[73,0] = 44533; $data[73,1] = 13; $data[73,2] = 25; $data[73,3] = 303.8405444822557
$data[74,0] = 44534; $data[74,1] = 1; $data[74,2] = 26; $data[74,3] = 315.994166261546
$data[75,0] = 44535; $data[75,1] = 2; $data[75,2] = 23; $data[75,3] = 279.5333009236753
$data[76,0] = 44536; $data[76,1] = 8; $data[76,2] = 30; $data[76,3] = 364.6086533787068
$data[77,0] = 44537; $data[77,1] = 16; $data[77,2] = 42; $data[77,3] = 510.4521147301896
$data[78,0] = 44538; $data[78,1] = 0; $data[78,2] = 41; $data[78,3] = 498.2984929508993

$ws.Range("A386:D464").Value = $data

Write-Host "Appended rows $firstRow to $lastRow ($numRows rows) to $($ws.Name)"
